$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 689.53352325
$schedule.Range("F2").Value = 11.40101724950397

# --- Sheet "Detailed" updates ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B17").Value = 1.77044
$detailed.Range("B18").Value = 0.7
$detailed.Range("C19").Value = "historical"
$detailed.Range("B20").Value = 0.7
$detailed.Range("C20").Value = "historical"
$detailed.Range("B21").Value = -0.94499
$detailed.Range("B22").Value = -2.54301
$detailed.Range("B23").Value = -5.50985
$detailed.Range("B24").Value = -5.58973
$detailed.Range("B25").Value = -0.97407
$detailed.Range("B26").Value = -6.23314
$detailed.Range("B27").Value = -5.50985
$detailed.Range("B28").Value = -2.83943
$detailed.Range("B29").Value = -2.74882
$detailed.Range("B30").Value = 0
$detailed.Range("B31").Value = 0.00047
$detailed.Range("B32").Value = 0.51
$detailed.Range("B35").Value = -2.49304
$detailed.Range("B36").Value = -0.00644
$detailed.Range("B37").Value = 3.10528
$detailed.Range("B38").Value = 13.52695
$detailed.Range("B39").Value = 43.33061
$detailed.Range("B40").Value = 56.98
$detailed.Range("B41").Value = 64.69326
$detailed.Range("B42").Value = 65
$detailed.Range("B43").Value = 67.48654000000001
$detailed.Range("B44").Value = 61.29106
$detailed.Range("B45").Value = 62.91375
$detailed.Range("B46").Value = 62.28206
$detailed.Range("B47").Value = 61.87129
$detailed.Range("B48").Value = 61.09939
$detailed.Range("B49").Value = 63.64693
